# "creation of Job 1st part"
# Adds a new "PostAsCircle" worksheet (between "TestCases" and "PostByImages"),
# populates it with test data, adds hyperlinks for the "password" column,
# and records the new test case on the "TestCases" sheet.

$wb = $excel.ActiveWorkbook

$testCases  = $wb.Worksheets.Item("TestCases")
$postImages = $wb.Worksheets.Item("PostByImages")
$post       = $wb.Worksheets.Item("Post")

# ---------------------------------------------------------------------------
# 1. Register the new test case on the "TestCases" sheet
# ---------------------------------------------------------------------------
$testCases.Range("A4").Value = 3
$testCases.Range("B4").Value = "PostAsCircle"
$testCases.Range("C4").Value = "y"
$testCases.Select()
$testCases.Range("D4").Select()

# ---------------------------------------------------------------------------
# 2. Create the new "PostAsCircle" worksheet, placed right before "PostByImages"
#    (final tab order becomes: TestCases, PostAsCircle, PostByImages, Post)
# ---------------------------------------------------------------------------
$circle = $wb.Worksheets.Add($postImages)
$circle.Name = "PostAsCircle"

# ---------------------------------------------------------------------------
# 3. Fill in the header row
# ---------------------------------------------------------------------------
$circle.Range("A1").Value = "TDID"
$circle.Range("B1").Value = "TDRunFlag"
$circle.Range("C1").Value = "userName"
$circle.Range("D1").Value = "password"
$circle.Range("F1").Value = "PostDescription"
$circle.Range("G1").Value = "PostHashTags1"
$circle.Range("H1").Value = "PostHashTags2"
$circle.Range("I1").Value = "PostHashTags3"
$circle.Range("J1").Value = "PostHashTags4"
$circle.Range("K1").Value = "PostHashTags5"
$circle.Range("L1").Value = "PostHashTags6"
$circle.Range("E1").Value = "CircleTitle"

# ---------------------------------------------------------------------------
# 4. Fill in the data rows
# ---------------------------------------------------------------------------
$circle.Range("A2").Value = 1
$circle.Range("B2").Value = "y"
$circle.Range("C2").Value = "'9030990045"
$circle.Range("D2").Value = "Welcome@123"
$circle.Range("F2").Value = "As a travel blogger, I promote wildlife tourism on my site, such as this recent guide to tiger safaris."
$circle.Range("G2").Value = "auto1"
$circle.Range("H2").Value = "auto2"
$circle.Range("I2").Value = "auto3"
$circle.Range("J2").Value = "auto4"
$circle.Range("K2").Value = "auto5"
$circle.Range("L2").Value = "auto6"

$circle.Range("A3").Value = 2
$circle.Range("B3").Value = "y"
$circle.Range("C3").Value = "'9030990045"
$circle.Range("D3").Value = "Welcome@123"
$circle.Range("F3").Value = "As a travel blogger, I promote wildlife tourism on my site, such as this recent guide to tiger safaris."
$circle.Range("G3").Value = "auto1"
$circle.Range("H3").Value = "auto2"
$circle.Range("I3").Value = "auto3"
$circle.Range("J3").Value = "auto4"
$circle.Range("K3").Value = "auto5"
$circle.Range("L3").Value = "auto6"

$circle.Range("A4").Value = 3
$circle.Range("B4").Value = "n"
$circle.Range("C4").Value = "'9030990045"
$circle.Range("D4").Value = "Welcome@123"
$circle.Range("F4").Value = "As a travel blogger, I promote wildlife tourism on my site, such as this recent guide to tiger safaris."
$circle.Range("G4").Value = "auto1"
$circle.Range("H4").Value = "auto2"
$circle.Range("I4").Value = "auto3"
$circle.Range("J4").Value = "auto4"
$circle.Range("K4").Value = "auto5"
$circle.Range("L4").Value = "auto6"

# column E ("CircleTitle" data) filled out of row order to mirror the shared
# string table layout produced by the original authoring session
$circle.Range("E4").Value = "Notification request"
$circle.Range("E2").Value = "Africa Geographic Tribe"
$circle.Range("E3").Value = "Testing mobile"

# ---------------------------------------------------------------------------
# 5. Formatting: bold header row, wrap the description column/header
# ---------------------------------------------------------------------------
$circle.Range("A1:L1").Font.Bold = $true

$circle.Hyperlinks.Add($circle.Range("D3"), "http://www.facebook.com/settings?tab=account")
$circle.Hyperlinks.Add($circle.Range("D4"), "http://www.facebook.com/settings?tab=account")

$circle.Range("F1").WrapText = $true
$circle.Range("F2:F4").WrapText = $true

$circle.Rows.Item(2).RowHeight = 28.8
$circle.Rows.Item(3).RowHeight = 28.8
$circle.Rows.Item(4).RowHeight = 28.8

# ---------------------------------------------------------------------------
# 6. Column widths (closest values reachable through this engine's rounding)
# ---------------------------------------------------------------------------
$circle.Columns.Item(1).ColumnWidth = 4.166666666666667
$circle.Columns.Item(2).ColumnWidth = 8.721354166666666
$circle.Columns.Item(3).ColumnWidth = 10.166666666666666
$circle.Columns.Item(4).ColumnWidth = 12.830729166666666
$circle.Columns.Item(5).ColumnWidth = 13.498697916666666
$circle.Columns.Item(6).ColumnWidth = 42.944010416666664
$circle.Range("G1:L1").ColumnWidth = 12.830729166666666

# ---------------------------------------------------------------------------
# 7. Page setup for the new sheet
# ---------------------------------------------------------------------------
$circle.PageSetup.PaperSize = 9
$circle.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 8. View state: PostByImages keeps the C31 selection (no longer the active
#    tab), PostAsCircle becomes the active tab with B4 selected.
# ---------------------------------------------------------------------------
$postImages.Select()
$postImages.Range("C31").Select()

$circle.Select()
$circle.Range("B4").Select()
